$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.269.88"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.594.03"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0604"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("D12").Value = "1.819.06"
$ws.Range("D13").Value = "1.580.58"
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "26.269.77"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.86%  "
$ws.Range("E19").Value = "  +3.86%  "
$ws.Range("D20").Value = "0.0₃0720"
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.25%  "
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").Value = "1.471.71"
$ws.Range("E32").Value = "  +4.26%  "
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.568"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.73%  "
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.74"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.25%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E42").Value = "  +1.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.934"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.01%  "
$ws.Range("D44").Value = "1.731.68"
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("E45").Value = "  -1.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.01%  "
$ws.Range("E48").Value = "  -0.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0501"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("E51").Value = "  +0.02%  "
